$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is purely numeric-looking text (e.g. "584.91").
# Force text storage (matching the source inlineStr/text cells) by setting
# NumberFormat to text ("@") before assigning, then resetting the cell
# style back to Normal so no stray style index is left on the cell.
$numericTextCells = @("D5", "D6", "D8", "D13", "D19", "D21", "D23", "D24", "D28", "D30", "D34", "D36", "D37", "D39", "D41", "D43", "D45", "D46", "D47", "D50", "D51")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '63.203.26'
$ws.Range("D3").Value = '2.569.59'
$ws.Range("E3").Value = '  +1.06%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '584.91'
$ws.Range("E5").Value = '  +3.05%  '
$ws.Range("D6").Value = '147.62'
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '0.602'
$ws.Range("E8").Value = '  +3.02%  '
$ws.Range("E9").Value = '  +3.74%  '
$ws.Range("E10").Value = '  +0.79%  '
$ws.Range("E11").Value = '  +0.11%  '
$ws.Range("E12").Value = '  +1.55%  '
$ws.Range("D13").Value = '27.51'
$ws.Range("E13").Value = '  +1.14%  '
$ws.Range("D14").Value = '3.029.93'
$ws.Range("E14").Value = '  +1.01%  '
$ws.Range("D15").Value = '63.150.85'
$ws.Range("E15").Value = '  +0.51%  '
$ws.Range("E16").Value = '  +3.84%  '
$ws.Range("D17").Value = '2.567.72'
$ws.Range("E17").Value = '  +0.97%  '
$ws.Range("E18").Value = '  -0.65%  '
$ws.Range("D19").Value = '343.37'
$ws.Range("E19").Value = '  +2.07%  '
$ws.Range("E20").Value = '  +3.33%  '
$ws.Range("D21").Value = '6.90'
$ws.Range("E21").Value = '  +2.49%  '
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").Value = '5.55'
$ws.Range("E23").Value = '  -3.60%  '
$ws.Range("D24").Value = '66.90'
$ws.Range("D25").Value = '2.697.19'
$ws.Range("E25").Value = '  +1.08%  '
$ws.Range("E26").Value = '  +1.23%  '
$ws.Range("E27").Value = '  +0.78%  '
$ws.Range("D28").Value = '8.17'
$ws.Range("E28").Value = '  +12.68%  '
$ws.Range("E29").Value = '  +1.67%  '
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("E31").Value = '  -1.21%  '
$ws.Range("E32").Value = '  +7.97%  '
$ws.Range("E33").Value = '  +2.36%  '
$ws.Range("D34").Value = '464.99'
$ws.Range("E34").Value = '  +13.48%  '
$ws.Range("E35").Value = '  +3.75%  '
$ws.Range("D36").Value = '176.09'
$ws.Range("E36").Value = '  -0.91%  '
$ws.Range("D37").Value = '0.408'
$ws.Range("E37").Value = '  +2.26%  '
$ws.Range("D39").Value = '4.58'
$ws.Range("E39").Value = '  +5.21%  '
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("D41").Value = '1.76'
$ws.Range("E41").Value = '  +0.12%  '
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("D43").Value = '152.11'
$ws.Range("E43").Value = '  -0.62%  '
$ws.Range("E44").Value = '  +2.34%  '
$ws.Range("D45").Value = '21.01'
$ws.Range("E45").Value = '  +1.09%  '
$ws.Range("D46").Value = '0.0548'
$ws.Range("E46").Value = '  +5.78%  '
$ws.Range("D47").Value = '0.613'
$ws.Range("E47").Value = '  +1.35%  '
$ws.Range("E48").Value = '  +2.24%  '
$ws.Range("E49").Value = '  +1.82%  '
$ws.Range("D50").Value = '1.75'
$ws.Range("E50").Value = '  -1.44%  '
$ws.Range("D51").Value = '11.39'
$ws.Range("E51").Value = '  -0.03%  '

foreach ($addr in $numericTextCells) {
    $ws.Range($addr).Style = "Normal"
}
